$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a text value into a cell without leaving it as a Number type,
# and without leaving a residual NumberFormat style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "63.703.61"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "3.292.64"
$ws.Range("E3").Value = "  +4.97%  "

$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue $ws.Range("D5") "600.49"
$ws.Range("E5").Value = "  +2.45%  "

Set-TextValue $ws.Range("D6") "141.42"
$ws.Range("E6").Value = "  +3.03%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.290.67"
$ws.Range("E8").Value = "  +4.96%  "

Set-TextValue $ws.Range("D9") "0.519"
$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("E10").Value = "  +2.92%  "

Set-TextValue $ws.Range("D11") "5.42"
$ws.Range("E11").Value = "  +3.66%  "

$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("E13").Value = "  +0.92%  "

Set-TextValue $ws.Range("D14") "34.44"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").Value = "3.837.25"
$ws.Range("E15").Value = "  +5.02%  "

$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").Value = "3.291.17"
$ws.Range("E17").Value = "  +5.07%  "

$ws.Range("D18").Value = "63.774.08"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("E19").Value = "  +2.49%  "

Set-TextValue $ws.Range("D20") "478.57"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("E21").Value = "  -0.68%  "

Set-TextValue $ws.Range("D22") "0.728"
$ws.Range("E22").Value = "  +4.45%  "

Set-TextValue $ws.Range("D23") "8.02"
$ws.Range("E23").Value = "  +4.81%  "

Set-TextValue $ws.Range("D24") "13.64"
$ws.Range("E24").Value = "  +5.35%  "

Set-TextValue $ws.Range("D25") "84.08"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("E26").Value = "  +0.11%  "

Set-TextValue $ws.Range("D27") "2.75"
$ws.Range("E27").Value = "  +2.17%  "

$ws.Range("E28").Value = "  +0.03%  "

Set-TextValue $ws.Range("D29") "7.23"
$ws.Range("E29").Value = "  +4.28%  "

Set-TextValue $ws.Range("D30") "8.06"
$ws.Range("E30").Value = "  +1.74%  "

Set-TextValue $ws.Range("D31") "2.15"
$ws.Range("E31").Value = "  +1.95%  "

Set-TextValue $ws.Range("D32") "28.56"
$ws.Range("E32").Value = "  +7.01%  "

Set-TextValue $ws.Range("D33") "0.105"
$ws.Range("E33").Value = "  -1.63%  "

Set-TextValue $ws.Range("D34") "2.53"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("E35").Value = "  +3.41%  "

Set-TextValue $ws.Range("D36") "5.96"
$ws.Range("E36").Value = "  +3.41%  "

Set-TextValue $ws.Range("D37") "53.18"
$ws.Range("E37").Value = "  +1.69%  "

$ws.Range("E38").Value = "  +6.11%  "

Set-TextValue $ws.Range("D39") "0.0397"
$ws.Range("E39").Value = "  +3.16%  "

Set-TextValue $ws.Range("D40") "425.99"
$ws.Range("E40").Value = "  +2.43%  "

$ws.Range("D41").Value = "3.064.77"
$ws.Range("E41").Value = "  +5.64%  "

Set-TextValue $ws.Range("D42") "8.32"
$ws.Range("E42").Value = "  +1.68%  "

Set-TextValue $ws.Range("D43") "2.72"
$ws.Range("E43").Value = "  +0.28%  "

Set-TextValue $ws.Range("D44") "0.113"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("E45").Value = "  +1.11%  "

Set-TextValue $ws.Range("D46") "2.18"
$ws.Range("E46").Value = "  +3.48%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "26.17"
$ws.Range("E47").Value = "  +3.16%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D48") "0.999"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "126.12"
$ws.Range("E49").Value = "  +4.47%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D50") "0.114"
$ws.Range("E50").Value = "  +1.66%  "

Set-TextValue $ws.Range("D51") "35.06"
$ws.Range("E51").Value = "  +12.01%  "
